$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 102 values (columns B..Z) with the revised figures
$ws.Range("B102").Value = 39368
$ws.Range("C102").Value = 29624
$ws.Range("D102").Value = 25311
$ws.Range("E102").Value = 2871
$ws.Range("F102").Value = 10877
$ws.Range("G102").Value = 11599
$ws.Range("H102").Value = 4486
$ws.Range("I102").Value = 7787
$ws.Range("J102").Value = 4641
$ws.Range("K102").Value = 3112
$ws.Range("L102").Value = 10937
$ws.Range("M102").Value = 10058
$ws.Range("N102").Value = 1725
$ws.Range("O102").Value = 4869
$ws.Range("P102").Value = 4433
$ws.Range("Q102").Value = 431
$ws.Range("R102").Value = 3392
$ws.Range("S102").Value = 853
$ws.Range("T102").Value = 12265
$ws.Range("U102").Value = 11356
$ws.Range("V102").Value = 252
$ws.Range("W102").Value = 1308
$ws.Range("X102").Value = 9752
$ws.Range("Y102").Value = 998
$ws.Range("Z102").Value = 38076

# Add new row 103 with the next quarterly period
# Force the period label to be stored as text (not auto-converted to a date
# serial) by temporarily applying a text number format, then restoring the
# default "Normal" style so no residual formatting is left on the cell.
$ws.Range("A103").NumberFormat = "@"
$ws.Range("A103").Value = "01-04-2021"
$ws.Range("A103").Style = "Normal"
$ws.Range("B103").Value = 41667
$ws.Range("C103").Value = 32833
$ws.Range("D103").Value = 26556
$ws.Range("E103").Value = 3547
$ws.Range("F103").Value = 10878
$ws.Range("G103").Value = 12230
$ws.Range("H103").Value = 6197
$ws.Range("I103").Value = 8218
$ws.Range("J103").Value = 4998
$ws.Range("K103").Value = 3193
$ws.Range("L103").Value = 10700
$ws.Range("M103").Value = 9811
$ws.Range("N103").Value = 762
$ws.Range("O103").Value = 5401
$ws.Range("P103").Value = 4898
$ws.Range("Q103").Value = 499
$ws.Range("R103").Value = 3588
$ws.Range("S103").Value = 866
$ws.Range("T103").Value = 12671
$ws.Range("U103").Value = 11804
$ws.Range("V103").Value = 271
$ws.Range("W103").Value = 1368
$ws.Range("X103").Value = 10124
$ws.Range("Y103").Value = 970
$ws.Range("Z103").Value = 39677
